$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the header label above the "current" ordnance column (R1): D5.1 -> D5.2
$ws.Range("R1").Value2 = "D5.2"

# Row 15 (AGM-154C): add a D5.2 expenditure of 8
$ws.Range("M15").Value2 = 8

# Row 15's remaining-stock cell (R15) now drops into the "low stock" warning
# styling (matching the look already used on R12/R14), achieved by copying
# the cell format from R12 which already carries that style.
$ws.Range("R12").Copy() | Out-Null
$ws.Range("R15").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

# Row 18 (GBU-16): add a D5.2 expenditure of 3
$ws.Range("M18").Value2 = 3

# Restore the current selection to R15, matching the saved view state
$ws.Activate()
$ws.Range("R15").Select()
